$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved for numeric-looking strings (prices, percentages)
# so Excel does not silently convert them to numbers (which would drop formatting
# like trailing zeros, e.g. "1.70" -> 1.7, or reformat large numbers with "." as thousand separators).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '43.080.97'
$ws.Range('E2').Value = '  -6.78%  '
$ws.Range('D3').Value = '2.551.16'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '298.99'
$ws.Range('E5').Value = '  -4.35%  '
$ws.Range('D6').Value = '94.26'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').Value = '0.575'
$ws.Range('E7').Value = '  -4.11%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -5.79%  '
$ws.Range('E10').Value = '  -8.08%  '
$ws.Range('D11').Value = '0.0809'
$ws.Range('E11').Value = '  -4.49%  '
$ws.Range('E12').Value = '  -5.25%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').Value = '2.943.68'
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('D15').Value = '2.556.50'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = '0.870'
$ws.Range('E16').Value = '  -5.68%  '
$ws.Range('D17').Value = '14.22'
$ws.Range('E17').Value = '  -4.61%  '
$ws.Range('D18').Value = '43.095.41'
$ws.Range('E18').Value = '  -7.19%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0981'
$ws.Range('E19').Value = '  -3.96%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '6.67'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').Value = '72.01'
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').Value = '261.84'
$ws.Range('E23').Value = '  -10.18%  '
$ws.Range('E24').Value = '  -4.89%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '2.15'
$ws.Range('E25').Value = '  -3.82%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '29.65'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '37.02'
$ws.Range('E29').Value = '  -6.46%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.13'
$ws.Range('E30').Value = '  -4.14%  '
$ws.Range('D31').Value = '5.97'
$ws.Range('E31').Value = '  -5.34%  '
$ws.Range('D32').Value = '154.54'
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').Value = '2.18'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').Value = '3.41'
$ws.Range('E34').Value = '  -5.71%  '
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('D36').Value = '0.0797'
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('E37').Value = '  -6.07%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -3.65%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '24.05'
$ws.Range('E39').Value = '  +13.20%  '
$ws.Range('D40').Value = '16.64'
$ws.Range('E40').Value = '  +5.60%  '
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('D42').Value = '0.0313'
$ws.Range('E42').Value = '  -6.18%  '
$ws.Range('D43').Value = '3.86'
$ws.Range('E43').Value = '  -4.40%  '
$ws.Range('D44').Value = '2.087.42'
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '86.04'
$ws.Range('E46').Value = '  -12.08%  '
$ws.Range('E47').Value = '  +3.56%  '
$ws.Range('D48').Value = '2.799.21'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '104.63'
$ws.Range('E49').Value = '  -4.30%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.70'
$ws.Range('E50').Value = '  -2.66%  '
$ws.Range('D51').Value = '8.72'
$ws.Range('E51').Value = '  -8.21%  '
